$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

$ws.Activate()
$ws.Range("R13").Select()
